$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells M1:O1, copying the formatting used by the existing header row (style of L1)
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("M1").Value = "microsoft_dlm"
$ws.Range("N1").Value = "recall_microsoft_dlm"
$ws.Range("O1").Value = "precision_microsoft_dlm"

# Data rows 2..62 for columns M (microsoft_dlm text), N (recall_microsoft_dlm), O (precision_microsoft_dlm)
$ws.Cells.Item(2, 13).Value = 'J''aimerais transférer 3564,00$ de mon compte chèque à mon compte d''épargne. J''aimerais verser 199$ à début de chaque mois sur mon compte Bell.'
$ws.Cells.Item(2, 14).Value = 1
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(3, 13).Value = 'Transfert 6000,07$ 80 et 9971$ de mon compte, celi le 11 novembre.'
$ws.Cells.Item(3, 14).Value = 0.8461538461538461
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(4, 13).Value = 'Transfert de 51766,00$ 6652$ de mon compte celi le 29 décembre une fois cette année.'
$ws.Cells.Item(4, 14).Value = 1
$ws.Cells.Item(4, 15).Value = 1
$ws.Cells.Item(5, 13).Value = 'Transfert 2013$ à mon compte air au 1 janvier 2022.'
$ws.Cells.Item(5, 14).Value = 1
$ws.Cells.Item(5, 15).Value = 0.875
$ws.Cells.Item(6, 13).Value = 'Je veux faire un dépôt de 5243,00$ aussi 37000,04$ 65 le 18 novembre 2021.'
$ws.Cells.Item(6, 14).Value = 0.9523809523809523
$ws.Cells.Item(6, 15).Value = 1
$ws.Cells.Item(7, 13).Value = 'Je dois faire un transfert international de 3666$ à mon oncle à Québec, le 31 octobre 2021.'
$ws.Cells.Item(7, 14).Value = 1
$ws.Cells.Item(7, 15).Value = 1
$ws.Cells.Item(8, 13).Value = 'J''aimerais, c''est du lait. Un paiement mensuel de 27,55$ sous qui sera versé tous les 15 du mois à partir d''octobre.'
$ws.Cells.Item(8, 14).Value = 1
$ws.Cells.Item(8, 15).Value = 1
$ws.Cells.Item(9, 13).Value = 'J''aimerais, c''est du lait. Un paiement mensuel de 65,33$ sous qui sera versé tous les 29 du mois à partir de mars 2022.'
$ws.Cells.Item(9, 14).Value = 1
$ws.Cells.Item(9, 15).Value = 1
$ws.Cells.Item(10, 13).Value = 'Transfert, 5000, quatre-vingts dollars à mon compte réer le 16 avril 2022.'
$ws.Cells.Item(10, 14).Value = 0.8888888888888888
$ws.Cells.Item(10, 15).Value = 1
$ws.Cells.Item(11, 13).Value = 'Moi j''aimerais transférer 2077$ au à Trois-Rivières s''il vous plaît de mon compte américain le 18 décembre.'
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = 1
$ws.Cells.Item(12, 13).Value = 'Transfert 33112,00$, 33214$.'
$ws.Cells.Item(12, 14).Value = 1
$ws.Cells.Item(12, 15).Value = 1
$ws.Cells.Item(13, 13).Value = 'J''aimerais transférer 4654$ de mon compte chèque à mon compte d''épargne, j''aimerais aussi verser 398$ à la fin de chaque mois sur mon compte Bell.'
$ws.Cells.Item(13, 14).Value = 1
$ws.Cells.Item(13, 15).Value = 1
$ws.Cells.Item(14, 13).Value = 'Transfère 8000,08$ 70 à mon contraire et à un autre 9071,56$ sous de mon quantité celi le 11 novembre.'
$ws.Cells.Item(14, 14).Value = 0.9285714285714286
$ws.Cells.Item(14, 15).Value = 0.9285714285714286
$ws.Cells.Item(15, 13).Value = 'Transfert de 31910,00$ 4452$ de mon compte celi le 23 décembre une fois cette année.'
$ws.Cells.Item(15, 14).Value = 1
$ws.Cells.Item(15, 15).Value = 1
$ws.Cells.Item(16, 13).Value = 'Transfert 5012$ à mon compte réer le 1 février 2022.'
$ws.Cells.Item(16, 14).Value = 1
$ws.Cells.Item(16, 15).Value = 0.875
$ws.Cells.Item(17, 13).Value = 'Je vais faire un dépôt de 90243,01$ autre 39000,03$ 61, le 23 novembre 2021.'
$ws.Cells.Item(17, 14).Value = 0.9583333333333334
$ws.Cells.Item(17, 15).Value = 1
$ws.Cells.Item(18, 13).Value = 'Je dois effectuer un transfert international de 8606$ à mon oncle à Paris le 30 juin 2021.'
$ws.Cells.Item(18, 14).Value = 1
$ws.Cells.Item(18, 15).Value = 1
$ws.Cells.Item(19, 13).Value = 'J''aimerais du lait. Un paiement mensuel de 39,68$ sous qui sera versé tous les 28 du mois à partir d''octobre.'
$ws.Cells.Item(19, 14).Value = 1
$ws.Cells.Item(19, 15).Value = 1
$ws.Cells.Item(20, 13).Value = 'J''aimerais céduler un paiement mensuel de 69,59$ sous qui sera versé tous les 14 du mois à partir de mars 2022.'
$ws.Cells.Item(20, 14).Value = 1
$ws.Cells.Item(20, 15).Value = 1
$ws.Cells.Item(21, 13).Value = 'Transférer 31000 et 4$ à mon compte réer au 16 avril 2022.'
$ws.Cells.Item(21, 14).Value = 1
$ws.Cells.Item(21, 15).Value = 1
$ws.Cells.Item(22, 13).Value = 'Moi, j''aimerais transférer 2078$ s''il vous plaît de mon compte américain le 17 décembre.'
$ws.Cells.Item(22, 14).Value = 1
$ws.Cells.Item(22, 15).Value = 1
$ws.Cells.Item(23, 13).Value = 'Transférer, 35130,00$ aussi 33920$.'
$ws.Cells.Item(23, 14).Value = 1
$ws.Cells.Item(23, 15).Value = 1
$ws.Cells.Item(24, 13).Value = 'Transférer 67030,00$ transférer aussi 25715,30$ sous.'
$ws.Cells.Item(24, 14).Value = 1
$ws.Cells.Item(24, 15).Value = 1
$ws.Cells.Item(25, 13).Value = 'J''aimerais transférer 2834,00$ de mon compte chèque à mon compte d''épargne et j''aimerais aussi verser 295$ à la fin de chaque mois sur mon compte belle.'
$ws.Cells.Item(25, 14).Value = 1
$ws.Cells.Item(25, 15).Value = 1
$ws.Cells.Item(26, 13).Value = 'Transfert 1000,00$ 10 et un autre 2072,12$ sous de mon compte celi le 20 novembre.'
$ws.Cells.Item(26, 14).Value = 0.9090909090909091
$ws.Cells.Item(26, 15).Value = 1
$ws.Cells.Item(27, 13).Value = 'Transfert 21110,03$ 1455,00$ de mon compte celi le 2 décembre de cette année.'
$ws.Cells.Item(27, 14).Value = 1
$ws.Cells.Item(27, 15).Value = 1
$ws.Cells.Item(28, 13).Value = 'Transfert, 13013$ à mon compte réel le 10 février 2022.'
$ws.Cells.Item(28, 14).Value = 1
$ws.Cells.Item(28, 15).Value = 1
$ws.Cells.Item(29, 13).Value = 'Je vais faire un dépôt de 81000 et 83,01$ autre transfert de 31000,09$ 160 et 3 le 30 novembre 2021.'
$ws.Cells.Item(29, 14).Value = 1
$ws.Cells.Item(29, 15).Value = 1
$ws.Cells.Item(30, 13).Value = 'Je dois effectuer un transfert Inter international de 9107$ à mon oncle le 30 juillet 2021.'
$ws.Cells.Item(30, 14).Value = 1
$ws.Cells.Item(30, 15).Value = 1
$ws.Cells.Item(31, 13).Value = 'J''aimerais céduler un paiement mensuel de 32,18$ sous le 7 de chaque mois à partir d''octobre.'
$ws.Cells.Item(31, 14).Value = 1
$ws.Cells.Item(31, 15).Value = 1
$ws.Cells.Item(32, 13).Value = 'Je m''assis du lait un paiement mensuel de 49,49$ sous tous les 11 du mois à partir d''avril 2022.'
$ws.Cells.Item(32, 14).Value = 1
$ws.Cells.Item(32, 15).Value = 0.9
$ws.Cells.Item(33, 13).Value = 'Transfert 21000 et 99$ à mon compte réer au 3 avril 2022.'
$ws.Cells.Item(33, 14).Value = 1
$ws.Cells.Item(33, 15).Value = 1
$ws.Cells.Item(34, 13).Value = 'Moi, j''aimerais transférer 20079$ en Floride. S''il vous plaît de mon compte américain, le 31 décembre.'
$ws.Cells.Item(34, 14).Value = 1
$ws.Cells.Item(34, 15).Value = 1
$ws.Cells.Item(35, 13).Value = 'Transférer 6103,00$ 34410$.'
$ws.Cells.Item(35, 14).Value = 1
$ws.Cells.Item(35, 15).Value = 1
$ws.Cells.Item(36, 13).Value = 'Je voudrais transférer 1345,00$ aussi 2564$ à mon compte réer aujourd?'
$ws.Cells.Item(36, 14).Value = 1
$ws.Cells.Item(36, 15).Value = 1
$ws.Cells.Item(37, 13).Value = 'Je voudrais transférer 2898,00$ 1135$ sur mon compte, cpg demain.'
$ws.Cells.Item(37, 14).Value = 1
$ws.Cells.Item(37, 15).Value = 1
$ws.Cells.Item(38, 13).Value = 'Je voudrais changer mon nip à 66 77 88.'
$ws.Cells.Item(38, 14).Value = 1
$ws.Cells.Item(38, 15).Value = 1
$ws.Cells.Item(39, 13).Value = 'Je voudrais faire 3 transferts sur mon compte d''épargne BNC. Le premier montant est de 1234$, le 2nd de 2345,00$, le dernier de 3456,00$.'
$ws.Cells.Item(39, 14).Value = 1
$ws.Cells.Item(39, 15).Value = 1
$ws.Cells.Item(40, 13).Value = 'Je voudrais transférer 6677,00$ 7788,00$ sur mon compte cpg demain.'
$ws.Cells.Item(40, 14).Value = 1
$ws.Cells.Item(40, 15).Value = 1
$ws.Cells.Item(41, 13).Value = 'Je voudrais transférer 1001,00$ 2002$ sur mon compte, CPG demain.'
$ws.Cells.Item(41, 14).Value = 1
$ws.Cells.Item(41, 15).Value = 1
$ws.Cells.Item(42, 13).Value = 'Je voudrais transférer 20001,00$ 3005$ sur mon compte, CPG demain.'
$ws.Cells.Item(42, 14).Value = 1
$ws.Cells.Item(42, 15).Value = 1
$ws.Cells.Item(43, 13).Value = 'Je voudrais transférer 8888,07$ 1567$ sur mon compte, cpg demain.'
$ws.Cells.Item(43, 14).Value = 1
$ws.Cells.Item(43, 15).Value = 1
$ws.Cells.Item(44, 13).Value = 'J''aimerais, c''est du lait. Un paiement mensuel de 35,25$ sous qui sera versé tous les 17 du mois à partir d''octobre.'
$ws.Cells.Item(44, 14).Value = 1
$ws.Cells.Item(44, 15).Value = 1
$ws.Cells.Item(45, 13).Value = 'J''aimerais du lait. Un paiement mensuel de 68,35$ sous qui traversaient tous les 8 du mois à partir d''octobre.'
$ws.Cells.Item(45, 14).Value = 1
$ws.Cells.Item(45, 15).Value = 1
$ws.Cells.Item(46, 13).Value = 'J''aimerais céduler un paiement mensuel de 58,09$ sous qui sera versé tous les premiers du mois à partir d''octobre.'
$ws.Cells.Item(46, 14).Value = 1
$ws.Cells.Item(46, 15).Value = 1
$ws.Cells.Item(47, 13).Value = 'Je voudrais faire 3 transferts sur mon compte d''épargne BNC. Le premier montant est de soi 65343$. Le 2nd de 2999,00$, le dernier de 20,00$.'
$ws.Cells.Item(47, 14).Value = 1
$ws.Cells.Item(47, 15).Value = 1
$ws.Cells.Item(48, 13).Value = 'Je voudrais faire 3 transferts sur mon compte d''épargne BNC. Le premier montant est de 1234$, le 2nd de 2345,00$, le dernier de 3456,00$.'
$ws.Cells.Item(48, 14).Value = 1
$ws.Cells.Item(48, 15).Value = 1
$ws.Cells.Item(49, 13).Value = 'Je voudrais faire 5 transferts sur mon compte d''épargne BNC, le premier d''un montant de 10$, le 2nd de 20$ le 3e de 30,00$ le 4e de 40,00$ le 5e de 50$.'
$ws.Cells.Item(49, 14).Value = 1
$ws.Cells.Item(49, 15).Value = 1
$ws.Cells.Item(50, 13).Value = 'Je voudrais faire 5 transferts sur mon compte d''épargne BNC, le premier de 60$, le 2nd de 70,00$, le 3e de quatre-vingts dollars, le 4e de 90,00$, le 5e de 100$.'
$ws.Cells.Item(50, 14).Value = 1
$ws.Cells.Item(50, 15).Value = 1
$ws.Cells.Item(51, 13).Value = 'Je voudrais faire 5 transferts sur mon compte d''épargne BNC, le premier de 100,00$, le 2e de 200$, le 3e de 300$, le 4e de 400,00$, le 5e de 500$.'
$ws.Cells.Item(51, 14).Value = 1
$ws.Cells.Item(51, 15).Value = 1
$ws.Cells.Item(52, 13).Value = 'Je voudrais faire 5 virements sur mon compte d''épargne BC le premier de 600$, le 2e de 700$, le 3e de 800,00$, le 4e de 900,00$, le 5e de 1000 Do.'
$ws.Cells.Item(52, 14).Value = 1
$ws.Cells.Item(52, 15).Value = 1
$ws.Cells.Item(53, 13).Value = 'Je voudrais faire 5 transferts sur mon compte d''épargne BNC, le premier de 1001,00$, le 2e de 1010$, le 3e de 1100,00$, le 4e de 1110,00$, le 5e de 1111$.'
$ws.Cells.Item(53, 14).Value = 1
$ws.Cells.Item(53, 15).Value = 1
$ws.Cells.Item(54, 13).Value = 'Je voudrais faire 5 transferts sur mon compte d''épargne BNC, le premier pour 2001,00$ le 2e pour 2020$, le 3e pour 2200$, 4e pour 2220,00$, le 5e pour 2222$.'
$ws.Cells.Item(54, 14).Value = 1
$ws.Cells.Item(54, 15).Value = 1
$ws.Cells.Item(55, 13).Value = 'Je voudrais faire 5 transferts sur mon compte réer, le premier pour 3001,00$, le 2e pour 3030$, le 3e pour 3300$, le 4e pour 3330$, le 5e pour 3333$.'
$ws.Cells.Item(55, 14).Value = 1
$ws.Cells.Item(55, 15).Value = 1
$ws.Cells.Item(56, 13).Value = 'Je voudrais payer 13511$ à ma marge de crédit.'
$ws.Cells.Item(56, 14).Value = 1
$ws.Cells.Item(56, 15).Value = 1
$ws.Cells.Item(57, 13).Value = 'Je voudrais payer 13533$ à ma marge de crédit.'
$ws.Cells.Item(57, 14).Value = 1
$ws.Cells.Item(57, 15).Value = 1
$ws.Cells.Item(58, 13).Value = 'Je voudrais payer 13555$ à ma marge de crédit.'
$ws.Cells.Item(58, 14).Value = 1
$ws.Cells.Item(58, 15).Value = 1
$ws.Cells.Item(59, 13).Value = 'Je voudrais payer 13566$ à ma marge de crédit.'
$ws.Cells.Item(59, 14).Value = 1
$ws.Cells.Item(59, 15).Value = 1
$ws.Cells.Item(60, 13).Value = 'Je voudrais payer 13577$ à ma marge de crédit.'
$ws.Cells.Item(60, 14).Value = 1
$ws.Cells.Item(60, 15).Value = 1
$ws.Cells.Item(61, 13).Value = 'Je voudrais payer 13588$ à ma marge de crédit.'
$ws.Cells.Item(61, 14).Value = 1
$ws.Cells.Item(61, 15).Value = 1
$ws.Cells.Item(62, 13).Value = 'Je voudrais payer 13599$ à ma marge de crédit.'
$ws.Cells.Item(62, 14).Value = 1
$ws.Cells.Item(62, 15).Value = 1
